$d = $word.ActiveDocument

$replacements = @(
    @{old="99×90="; new="14×58="},
    @{old="25×13="; new="44×62="},
    @{old="13×43="; new="34×18="},
    @{old="80×13="; new="39×28="},
    @{old="81×79="; new="95×41="},
    @{old="56×17="; new="30×73="},
    @{old="22×94="; new="25×66="},
    @{old="95×25="; new="46×80="},
    @{old="93×27="; new="83×17="},
    @{old="99×38="; new="53×21="},
    @{old="60×53="; new="43×33="},
    @{old="30×59="; new="78×68="},
    @{old="39×16="; new="65×27="},
    @{old="64×30="; new="52×60="},
    @{old="71×43="; new="89×54="},
    @{old="20×16="; new="89×17="},
    @{old="48×15="; new="76×20="},
    @{old="18×12="; new="64×21="},
    @{old="32×41="; new="44×30="},
    @{old="16×80="; new="25×45="},
    @{old="68×14="; new="48×72="},
    @{old="89×16="; new="41×69="},
    @{old="46×78="; new="43×56="},
    @{old="37×54="; new="44×26="},
    @{old="70×55="; new="83×65="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
